$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.690820217132568
$ws.Range("B1").Value = 3.69589376449585
$ws.Range("C1").Value = 3.423209667205811
$ws.Range("D1").Value = 3.330772399902344
$ws.Range("E1").Value = 1.200621843338013
